$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.122.52"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "3.422.33"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.06"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.89"
$ws.Range("E6").Value = "  -3.73%  "

$ws.Range("E7").Value = "  +9.25%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.741"
$ws.Range("E9").Value = "  +7.85%  "

$ws.Range("E10").Value = "  +16.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.84"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000221"
$ws.Range("E12").Value = "  +68.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.17"
$ws.Range("E13").Value = "  +8.17%  "

$ws.Range("D15").Value = "3.958.76"
$ws.Range("E15").Value = "  -0.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.09"
$ws.Range("E16").Value = "  +5.82%  "

$ws.Range("D17").Value = "3.417.95"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("E18").Value = "  +6.42%  "

$ws.Range("E19").Value = "  +5.83%  "

$ws.Range("D20").Value = "62.138.74"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.02"
$ws.Range("E21").Value = "  +40.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.27"
$ws.Range("E22").Value = "  +8.19%  "

$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.16"
$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("E25").Value = "  +3.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.66"
$ws.Range("E26").Value = "  +12.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.92"
$ws.Range("E27").Value = "  +8.44%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.75"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.01"
$ws.Range("E31").Value = "  +5.63%  "

$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.06"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("E36").Value = "  +3.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.12"
$ws.Range("E37").Value = "  +5.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  +8.87%  "

$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("E42").Value = "  +1.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.40"
$ws.Range("E43").Value = "  +2.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.22"
$ws.Range("E44").Value = "  +4.33%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("E46").Value = "  +7.65%  "

$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.32"
$ws.Range("E48").Value = "  +4.35%  "

$ws.Range("D49").Value = "3.773.53"
$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").Value = "2.119.64"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.99"
$ws.Range("E51").Value = "  +25.43%  "
